$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.541.98"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.12"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "263.87"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5228"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3246"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06808"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.78"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7789"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07780"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.841.68"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.47"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.023"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.98"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007970"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.562.64"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.087.23"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.623"
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.456"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.003"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.07"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.170"
$ws.Range("E26").Value = "  -7.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.678"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.01"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.90"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.186"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08746"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.112"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04832"
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7230"
$ws.Range("E34").Value = "  +4.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.862"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.089"
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01794"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.218"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4856"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.37"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8890"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.046"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.636"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4208"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05893"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.068"
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1240"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.98"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8896"
$ws.Range("E51").Value = "  +4.23%  "
